# Applies the cryptos.xlsx price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to Text format so numeric-looking strings (e.g. "1.00")
    # are stored verbatim instead of being coerced into a Number,
    # then drop back to the Normal style so no stray number format sticks.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '27.626.37'
$ws.Range("D3").Value = '1.666.67'
$ws.Range("E3").Value = '  -3.83%  '
Set-TextValue "D5" '215.64'
$ws.Range("E5").Value = '  -2.30%  '
Set-TextValue "D6" '0.512'
$ws.Range("E6").Value = '  -2.31%  '
$ws.Range("E7").Value = '  +0.09%  '
Set-TextValue "D8" '24.18'
$ws.Range("E8").Value = '  -0.87%  '
Set-TextValue "D11" '0.0878'
$ws.Range("E11").Value = '  -2.26%  '
$ws.Range("D12").Value = '1.902.45'
$ws.Range("E12").Value = '  -3.75%  '
$ws.Range("D13").Value = '1.675.85'
$ws.Range("E13").Value = '  -3.29%  '
$ws.Range("E14").Value = '  -3.34%  '
Set-TextValue "D15" '0.567'
$ws.Range("E15").Value = '  +0.78%  '
Set-TextValue "D16" '66.47'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").Value = '27.610.58'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E19").Value = '  -3.40%  '
Set-TextValue "D20" '7.72'
$ws.Range("E20").Value = '  -4.06%  '
$ws.Range("E21").Value = '  +0.10%  '
Set-TextValue "D22" '4.51'
$ws.Range("E22").Value = '  -3.26%  '
Set-TextValue "D23" '9.38'
$ws.Range("E23").Value = '  -3.56%  '
Set-TextValue "D24" '2.05'
$ws.Range("E24").Value = '  -3.63%  '
$ws.Range("E26").Value = '  -4.11%  '
Set-TextValue "D27" '16.43'
$ws.Range("E27").Value = '  -1.95%  '
Set-TextValue "D28" '1.00'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").Value = '  -2.40%  '
$ws.Range("E30").Value = '  +1.92%  '
Set-TextValue "D31" '0.0505'
$ws.Range("E31").Value = '  -1.68%  '
$ws.Range("E32").Value = '  -2.72%  '
$ws.Range("D33").Value = '1.462.46'
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  -4.47%  '
$ws.Range("E35").Value = '  -4.97%  '
$ws.Range("E36").Value = '  -1.37%  '
Set-TextValue "D37" '0.929'
$ws.Range("E37").Value = '  -4.29%  '
$ws.Range("E38").Value = '  -4.95%  '
$ws.Range("E39").Value = '  -1.99%  '
Set-TextValue "D40" '69.69'
$ws.Range("E40").Value = '  -1.67%  '
$ws.Range("E41").Value = '  -4.05%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  -3.33%  '
$ws.Range("E44").Value = '  -5.80%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D45" '0.794'
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.810.16'
$ws.Range("E46").Value = '  -3.70%  '
$ws.Range("E47").Value = '  +0.27%  '
Set-TextValue "D48" '88.99'
$ws.Range("E49").Value = '  -6.11%  '
$ws.Range("E50").Value = '  -2.62%  '
Set-TextValue "D51" '7.92'
$ws.Range("E51").Value = '  -4.29%  '
